# Update "想去人数" (interested-count) figures in column F across the
# workbook's sheets, reflecting refreshed scrape data
# (commit: "Update gh-pages to output generated at 456a3b4").
#
# Sheet order (per xl/workbook.xml): 1=展览, 2=演出, 3=本地生活(empty), 4=全部类型

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)
$sheet1Updates = @{
    3  = 519
    4  = 1515
    5  = 151
    8  = 149
    9  = 734
    10 = 1045
    11 = 62
    12 = 327
    13 = 51
    14 = 6386
    15 = 5
    18 = 150
    20 = 15271
    21 = 1516
    22 = 279
    23 = 139
    24 = 100
    25 = 11031
    26 = 749
    27 = 4312
    28 = 234
    30 = 16
    31 = 301
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 6).Value = 343

# --- Sheet 4: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item(4)
$sheet4Updates = @{
    3  = 519
    4  = 1515
    5  = 151
    7  = 343
    9  = 149
    10 = 734
    12 = 1045
    13 = 62
    14 = 327
    15 = 51
    17 = 6386
    18 = 5
    21 = 150
    23 = 15271
    24 = 1516
    25 = 279
    26 = 139
    27 = 100
    28 = 11031
    29 = 749
    30 = 4312
    31 = 234
    33 = 16
    34 = 301
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
